$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("2025")
$ws1.Range("B2").Value = 9739.537847600008
$ws1.Range("E2").Value = 289823.7596598056
$ws1.Range("I2").Value = 161752.8135478
$ws1.Range("L2").Value = 485245.29503538
$ws1.Range("M2").Value = 105905.87968015
$ws1.Range("N2").Value = 70831.955579581
$ws1.Range("O2").Value = 69610.4422391004

$ws2 = $wb.Worksheets.Item("2030")
$ws2.Range("B2").Value = 47386.06393082884
$ws2.Range("E2").Value = 271236.7992183856
$ws2.Range("I2").Value = 280426.171173861
$ws2.Range("L2").Value = 184420.4799505123
$ws2.Range("M2").Value = 113936.92264746
$ws2.Range("N2").Value = 33931.8246116005
$ws2.Range("O2").Value = 50485.47232467777

$ws3 = $wb.Worksheets.Item("2035")
$ws3.Range("A2").Value = 28619.61401238371
$ws3.Range("B2").Value = 23143.29485244409
$ws3.Range("E2").Value = 111916.8406725409
$ws3.Range("I2").Value = 150385.2728707001
$ws3.Range("M2").Value = 34803.41203795493
$ws3.Range("N2").Value = 44938.11408779013
$ws3.Range("O2").Value = 26938.31306104351

$ws4 = $wb.Worksheets.Item("2040")
$ws4.Range("N2").Value = 1014.766490779938

$ws5 = $wb.Worksheets.Item("2045")
$ws5.Range("A2").Value = 34409.11717595647
$ws5.Range("N2").Value = 5182.698656944208
$ws5.Range("O2").Value = 22972.54525065906
